$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top, shifting existing data down
$ws.Rows.Item(1).Insert()

# Add header labels in the newly inserted row
$ws.Range("A1").Value = "NUTS1"
$ws.Range("B1").Value = "State"

# Move selection to A2, matching the post-edit selection state
$ws.Range("A2").Select()
